$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3632.3333
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 3632.3333
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 3632.3333
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -3982.3333
$ws.Range("H28").Value = 542.5
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H53").Value = 8069.923
$ws.Range("I53").Value = 375
$ws.Range("K53").Value = 375
$ws.Range("M53").Value = 262
$ws.Range("H62").Value = 4425.2144
$ws.Range("I62").Value = 4457.923
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 4457.923
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -3833.923
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 4425.2144
$ws.Range("I65").Value = 4457.923
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 22289.615
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -19169.615
$ws.Range("N65").Value = -26240
$ws.Range("H74").Value = 7703.2856
$ws.Range("I74").Value = 6990
$ws.Range("K74").Value = 6990
$ws.Range("M74").Value = -6054
$ws.Range("H76").Value = 2929751
$ws.Range("I76").Value = 4279653
$ws.Range("K76").Value = 4279653
$ws.Range("M76").Value = -4279338
$ws.Range("H77").Value = 7703.2856
$ws.Range("I77").Value = 6990
$ws.Range("K77").Value = 34950
$ws.Range("M77").Value = -30270
$ws.Range("H79").Value = 2929751
$ws.Range("I79").Value = 4279653
$ws.Range("K79").Value = 4279653
$ws.Range("M79").Value = -4278561
$ws.Range("H80").Value = 7606.6
$ws.Range("J80").Value = 9736.272000000001
$ws.Range("L80").Value = 29208.816
$ws.Range("N80").Value = -31204.816
$ws.Range("H83").Value = 7606.6
$ws.Range("J83").Value = 9736.272000000001
$ws.Range("L83").Value = 87626.448
$ws.Range("N83").Value = -97610.448
$ws.Range("H86").Value = 2518.25
$ws.Range("I86").Value = 3240.6667
$ws.Range("J86").Value = 2084.8
$ws.Range("K86").Value = 3240.6667
$ws.Range("L86").Value = 2084.8
$ws.Range("M86").Value = -2117.6667
$ws.Range("N86").Value = -4330.8
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H89").Value = 2518.25
$ws.Range("I89").Value = 3240.6667
$ws.Range("J89").Value = 2084.8
$ws.Range("K89").Value = 16203.3335
$ws.Range("L89").Value = 10424
$ws.Range("M89").Value = -10587.3335
$ws.Range("N89").Value = -21656
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H92").Value = 2387.35
$ws.Range("I92").Value = 1952.9375
$ws.Range("J92").Value = 4125
$ws.Range("K92").Value = 1952.9375
$ws.Range("L92").Value = 4125
$ws.Range("M92").Value = -704.9375
$ws.Range("N92").Value = -6621
$ws.Range("H98").Value = 2850.625
$ws.Range("I98").Value = 2687.5
$ws.Range("K98").Value = 2687.5
$ws.Range("M98").Value = -1189.5
$ws.Range("H100").Value = 6546.1816
$ws.Range("I100").Value = 5429
$ws.Range("J100").Value = 8501.25
$ws.Range("K100").Value = 5429
$ws.Range("L100").Value = 8501.25
$ws.Range("M100").Value = -4888
$ws.Range("N100").Value = -9583.25
$ws.Range("H107").Value = 27778474
$ws.Range("I107").Value = 27778474
$ws.Range("K107").Value = 27778474
$ws.Range("M107").Value = -27776554
$ws.Range("H111").Value = 12350060
$ws.Range("I111").Value = 22226562
$ws.Range("J111").Value = 4433
$ws.Range("K111").Value = 66679686
$ws.Range("L111").Value = 13299
$ws.Range("M111").Value = -66676619
$ws.Range("N111").Value = -19433
$ws.Range("H113").Value = 6946.5
$ws.Range("J113").Value = 7214.364
$ws.Range("L113").Value = 7214.364
$ws.Range("N113").Value = -13722.364
$ws.Range("H116").Value = 5466.1577
$ws.Range("J116").Value = 5513.5
$ws.Range("L116").Value = 5513.5
$ws.Range("N116").Value = -12397.5
$ws.Range("H122").Value = 2850.625
$ws.Range("I122").Value = 2687.5
$ws.Range("K122").Value = 8062.5
$ws.Range("M122").Value = -5612.5
$ws.Range("H132").Value = 19234450
$ws.Range("I132").Value = 23813180
$ws.Range("J132").Value = 3788.3
$ws.Range("K132").Value = 71439540
$ws.Range("L132").Value = 11364.9
$ws.Range("M132").Value = -71437010
$ws.Range("N132").Value = -16424.9
$ws.Range("H135").Value = 1002.8
$ws.Range("I135").Value = 897.4857
$ws.Range("J135").Value = 1740
$ws.Range("K135").Value = 8077.3713
$ws.Range("L135").Value = 15660
$ws.Range("M135").Value = -5542.3713
$ws.Range("N135").Value = -20730
$ws.Range("H137").Value = 142067.92
$ws.Range("I137").Value = 166625.73
$ws.Range("J137").Value = 7000
$ws.Range("K137").Value = 499877.1900000001
$ws.Range("L137").Value = 21000
$ws.Range("M137").Value = -497327.1900000001
$ws.Range("N137").Value = -26100
$ws.Range("H138").Value = 3354.2935
$ws.Range("I138").Value = 897.9524
$ws.Range("J138").Value = 4080.817
$ws.Range("K138").Value = 2693.8572
$ws.Range("L138").Value = 12242.451
$ws.Range("M138").Value = 2446.1428
$ws.Range("N138").Value = -22522.451
$ws.Range("H141").Value = 11249.733
$ws.Range("I141").Value = 5193.5415
$ws.Range("J141").Value = 35474.5
$ws.Range("K141").Value = 15580.6245
$ws.Range("L141").Value = 106423.5
$ws.Range("M141").Value = -10400.6245
$ws.Range("N141").Value = -116783.5

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 401
$ws.Range("I5").Value = 400
$ws.Range("J5").Value = 402
$ws.Range("K5").Value = 400
$ws.Range("L5").Value = 402
$ws.Range("M5").Value = -288
$ws.Range("N5").Value = -626
$ws.Range("H32").Value = 6491.8726
$ws.Range("I32").Value = 4097.3423
$ws.Range("J32").Value = 11844.353
$ws.Range("K32").Value = 4097.3423
$ws.Range("L32").Value = 11844.353
$ws.Range("M32").Value = -3810.3423
$ws.Range("N32").Value = -12418.353
$ws.Range("H45").Value = 7996181.5
$ws.Range("I45").Value = 14387306
$ws.Range("K45").Value = 14387306
$ws.Range("M45").Value = -14386929
$ws.Range("H56").Value = 14933
$ws.Range("J56").Value = 14933
$ws.Range("L56").Value = 14933
$ws.Range("N56").Value = -16417
$ws.Range("H61").Value = 4306.7896
$ws.Range("I61").Value = 5031.6
$ws.Range("J61").Value = 3501.4443
$ws.Range("K61").Value = 5031.6
$ws.Range("L61").Value = 3501.4443
$ws.Range("M61").Value = -4819.6
$ws.Range("N61").Value = -3925.4443
$ws.Range("H62").Value = 100001
$ws.Range("J62").Value = 100001
$ws.Range("L62").Value = 100001
$ws.Range("N62").Value = -101249
$ws.Range("H63").Value = 4921.6
$ws.Range("I63").Value = 4451.5
$ws.Range("J63").Value = 5235
$ws.Range("K63").Value = 4451.5
$ws.Range("L63").Value = 5235
$ws.Range("M63").Value = -3765.5
$ws.Range("N63").Value = -6607
$ws.Range("H65").Value = 100001
$ws.Range("J65").Value = 100001
$ws.Range("L65").Value = 300003
$ws.Range("N65").Value = -306243
$ws.Range("H66").Value = 4921.6
$ws.Range("I66").Value = 4451.5
$ws.Range("J66").Value = 5235
$ws.Range("K66").Value = 22257.5
$ws.Range("L66").Value = 26175
$ws.Range("M66").Value = -18825.5
$ws.Range("N66").Value = -33039
$ws.Range("H68").Value = 67000
$ws.Range("J68").Value = 98000.5
$ws.Range("L68").Value = 98000.5
$ws.Range("N68").Value = -99622.5
$ws.Range("H71").Value = 67000
$ws.Range("J71").Value = 98000.5
$ws.Range("L71").Value = 294001.5
$ws.Range("N71").Value = -302113.5
$ws.Range("H74").Value = 90911.10000000001
$ws.Range("I74").Value = 40304.78
$ws.Range("J74").Value = 257189
$ws.Range("K74").Value = 40304.78
$ws.Range("L74").Value = 257189
$ws.Range("M74").Value = -39430.78
$ws.Range("N74").Value = -258937
$ws.Range("H77").Value = 90911.10000000001
$ws.Range("I77").Value = 40304.78
$ws.Range("J77").Value = 257189
$ws.Range("K77").Value = 201523.9
$ws.Range("L77").Value = 1285945
$ws.Range("M77").Value = -197155.9
$ws.Range("N77").Value = -1294681
$ws.Range("H80").Value = 60000
$ws.Range("J80").Value = 60000
$ws.Range("L80").Value = 60000
$ws.Range("N80").Value = -61996
$ws.Range("H83").Value = 60000
$ws.Range("J83").Value = 60000
$ws.Range("L83").Value = 180000
$ws.Range("N83").Value = -189984
$ws.Range("H88").Value = 1153.3846
$ws.Range("I88").Value = 625
$ws.Range("K88").Value = 625
$ws.Range("M88").Value = -219
$ws.Range("H91").Value = 1153.3846
$ws.Range("I91").Value = 625
$ws.Range("K91").Value = 625
$ws.Range("M91").Value = 779
$ws.Range("H92").Value = 48550
$ws.Range("J92").Value = 48550
$ws.Range("L92").Value = 48550
$ws.Range("N92").Value = -53542
$ws.Range("H94").Value = 43331.668
$ws.Range("J94").Value = 43331.668
$ws.Range("L94").Value = 43331.668
$ws.Range("N94").Value = -45133.668
$ws.Range("H102").Value = 2780954.8
$ws.Range("I102").Value = 3969829.2
$ws.Range("K102").Value = 3969829.2
$ws.Range("M102").Value = -3968207.2
$ws.Range("H108").Value = 31895
$ws.Range("J108").Value = 31895
$ws.Range("L108").Value = 31895
$ws.Range("N108").Value = -39575
$ws.Range("H110").Value = 2138398.5
$ws.Range("I110").Value = 2526780
$ws.Range("K110").Value = 2526780
$ws.Range("M110").Value = -2524735
$ws.Range("H130").Value = 50107
$ws.Range("J130").Value = 50107
$ws.Range("L130").Value = 50107
$ws.Range("N130").Value = -60147
$ws.Range("H132").Value = 2394.0596
$ws.Range("I132").Value = 1803.78
$ws.Range("K132").Value = 5411.34
$ws.Range("M132").Value = -2881.34
$ws.Range("H134").Value = 79000
$ws.Range("J134").Value = 79000
$ws.Range("L134").Value = 79000
$ws.Range("N134").Value = -89140
$ws.Range("H135").Value = 200066670
$ws.Range("I135").Value = 23390
$ws.Range("J135").Value = 250077500
$ws.Range("K135").Value = 23390
$ws.Range("L135").Value = 250077500
$ws.Range("M135").Value = -18320
$ws.Range("N135").Value = -250087640
$ws.Range("H136").Value = 4306.7896
$ws.Range("I136").Value = 5031.6
$ws.Range("J136").Value = 3501.4443
$ws.Range("K136").Value = 15094.8
$ws.Range("L136").Value = 10504.3329
$ws.Range("M136").Value = -12544.8
$ws.Range("N136").Value = -15604.3329
$ws.Range("H141").Value = 64332.668
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 64332.668
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 64332.668
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -74692.66800000001

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 401
$ws.Range("I4").Value = 400
$ws.Range("J4").Value = 402
$ws.Range("K4").Value = 400
$ws.Range("L4").Value = 402
$ws.Range("M4").Value = -285
$ws.Range("N4").Value = -632
$ws.Range("H5").Value = 3132
$ws.Range("I5").Value = 297.25
$ws.Range("J5").Value = 5021.8335
$ws.Range("K5").Value = 297.25
$ws.Range("L5").Value = 5021.8335
$ws.Range("M5").Value = -184.25
$ws.Range("N5").Value = -5247.8335
$ws.Range("H22").Value = 879.4
$ws.Range("I22").Value = 974.25
$ws.Range("K22").Value = 974.25
$ws.Range("M22").Value = -801.25
$ws.Range("H29").Value = 120490.2
$ws.Range("I29").Value = 200272
$ws.Range("K29").Value = 200272
$ws.Range("M29").Value = -199983
$ws.Range("H68").Value = 52655
$ws.Range("J68").Value = 52655
$ws.Range("L68").Value = 52655
$ws.Range("N68").Value = -54277
$ws.Range("H69").Value = 49000.5
$ws.Range("J69").Value = 49000.5
$ws.Range("L69").Value = 49000.5
$ws.Range("N69").Value = -50622.5
$ws.Range("H71").Value = 52655
$ws.Range("J71").Value = 52655
$ws.Range("L71").Value = 157965
$ws.Range("N71").Value = -166077
$ws.Range("H72").Value = 49000.5
$ws.Range("J72").Value = 49000.5
$ws.Range("L72").Value = 147001.5
$ws.Range("N72").Value = -155113.5
$ws.Range("H86").Value = 4354977.5
$ws.Range("I86").Value = 4769309
$ws.Range("K86").Value = 4769309
$ws.Range("M86").Value = -4768186
$ws.Range("H89").Value = 4354977.5
$ws.Range("I89").Value = 4769309
$ws.Range("K89").Value = 23846545
$ws.Range("M89").Value = -23840929
$ws.Range("H94").Value = 2783764.8
$ws.Range("I94").Value = 3125696.2
$ws.Range("K94").Value = 3125696.2
$ws.Range("M94").Value = -3125245.2
$ws.Range("H99").Value = 79651.30499999999
$ws.Range("J99").Value = 4247.5
$ws.Range("L99").Value = 4247.5
$ws.Range("N99").Value = -7243.5
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H134").Value = 3701.125
$ws.Range("I134").Value = 1601.7778
$ws.Range("K134").Value = 4805.3334
$ws.Range("M134").Value = -2270.3334

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 82131
$ws.Range("J9").Value = 82131
$ws.Range("L9").Value = 82131
$ws.Range("N9").Value = -82467
$ws.Range("H16").Value = 2185.0715
$ws.Range("I16").Value = 1224.5
$ws.Range("K16").Value = 1224.5
$ws.Range("M16").Value = -937.5
$ws.Range("H22").Value = 470.75
$ws.Range("I22").Value = 354
$ws.Range("J22").Value = 587.5
$ws.Range("K22").Value = 354
$ws.Range("L22").Value = 587.5
$ws.Range("M22").Value = -4
$ws.Range("N22").Value = -1287.5
$ws.Range("H31").Value = 20111.826
$ws.Range("I31").Value = 1766.5
$ws.Range("K31").Value = 1766.5
$ws.Range("M31").Value = -1471.5
$ws.Range("H34").Value = 20111.826
$ws.Range("I34").Value = 1766.5
$ws.Range("K34").Value = 1766.5
$ws.Range("M34").Value = -1564.5
$ws.Range("H57").Value = 54200
$ws.Range("I57").Value = 54200
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 54200
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -53640
$ws.Range("N57").ClearContents()
$ws.Range("H58").Value = 1723.0667
$ws.Range("I58").Value = 1363.1305
$ws.Range("K58").Value = 1363.1305
$ws.Range("M58").Value = -1160.1305
$ws.Range("H82").Value = 59998.5
$ws.Range("J82").Value = 59998.5
$ws.Range("L82").Value = 59998.5
$ws.Range("N82").Value = -60720.5
$ws.Range("H85").Value = 59998.5
$ws.Range("J85").Value = 59998.5
$ws.Range("L85").Value = 59998.5
$ws.Range("N85").Value = -62494.5
$ws.Range("H87").Value = 40000
$ws.Range("J87").Value = 40000
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42372
$ws.Range("H88").Value = 41022.6
$ws.Range("J88").Value = 41022.6
$ws.Range("L88").Value = 41022.6
$ws.Range("N88").Value = -41834.6
$ws.Range("H90").Value = 40000
$ws.Range("J90").Value = 40000
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -131856
$ws.Range("H91").Value = 41022.6
$ws.Range("J91").Value = 41022.6
$ws.Range("L91").Value = 41022.6
$ws.Range("N91").Value = -43830.6
$ws.Range("H107").Value = 32259146
$ws.Range("I107").Value = 1117.4333
$ws.Range("K107").Value = 1117.4333
$ws.Range("M107").Value = 802.5667000000001
$ws.Range("H113").Value = 2185.0715
$ws.Range("I113").Value = 1224.5
$ws.Range("K113").Value = 1224.5
$ws.Range("M113").Value = 945.5
$ws.Range("H132").Value = 37528.73
$ws.Range("I132").Value = 2968.0833
$ws.Range("K132").Value = 8904.249899999999
$ws.Range("M132").Value = -6374.249899999999
$ws.Range("H134").Value = 40983.56
$ws.Range("I134").Value = 57556.117
$ws.Range("K134").Value = 172668.351
$ws.Range("M134").Value = -170133.351
$ws.Range("H136").Value = 1723.0667
$ws.Range("I136").Value = 1363.1305
$ws.Range("K136").Value = 4089.3915
$ws.Range("M136").Value = -1539.3915

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 56960
$ws.Range("J37").Value = 56960
$ws.Range("L37").Value = 170880
$ws.Range("N37").Value = -171104
$ws.Range("H51").Value = 3955.8667
$ws.Range("I51").Value = 534
$ws.Range("K51").Value = 1602
$ws.Range("M51").Value = -1142
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H107").Value = 526.55554
$ws.Range("I107").Value = 257.9091
$ws.Range("J107").Value = 711.25
$ws.Range("K107").Value = 773.7273
$ws.Range("L107").Value = 2133.75
$ws.Range("M107").Value = 1146.2727
$ws.Range("N107").Value = -5973.75
$ws.Range("H114").Value = 17544542
$ws.Range("I114").Value = 41666852
$ws.Range("J114").Value = 1045.091
$ws.Range("K114").Value = 125000556
$ws.Range("L114").Value = 3135.273
$ws.Range("M114").Value = -124997302
$ws.Range("N114").Value = -9643.272999999999
$ws.Range("H122").Value = 1076.1666
$ws.Range("I122").Value = 1062.1
$ws.Range("J122").Value = 1093.75
$ws.Range("K122").Value = 9558.9
$ws.Range("L122").Value = 9843.75
$ws.Range("M122").Value = -7108.9
$ws.Range("N122").Value = -14743.75
$ws.Range("H131").Value = 9472469
$ws.Range("I131").Value = 7576572.5
$ws.Range("J131").Value = 10104435
$ws.Range("K131").Value = 22729717.5
$ws.Range("L131").Value = 30313305
$ws.Range("M131").Value = -22724677.5
$ws.Range("N131").Value = -30323385
$ws.Range("H132").Value = 2192.853
$ws.Range("J132").Value = 2602.2307
$ws.Range("L132").Value = 23420.0763
$ws.Range("N132").Value = -28480.0763
$ws.Range("H133").Value = 4750

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 12178.909
$ws.Range("J12").Value = 2332.6667
$ws.Range("L12").Value = 2332.6667
$ws.Range("N12").Value = -2612.6667
$ws.Range("H80").Value = 43868228
$ws.Range("I80").Value = 74915200
$ws.Range("J80").Value = 402459.8
$ws.Range("K80").Value = 74915200
$ws.Range("L80").Value = 402459.8
$ws.Range("M80").Value = -74914202
$ws.Range("N80").Value = -404455.8
$ws.Range("H83").Value = 43868228
$ws.Range("I83").Value = 74915200
$ws.Range("J83").Value = 402459.8
$ws.Range("K83").Value = 374576000
$ws.Range("L83").Value = 2012299
$ws.Range("M83").Value = -374571008
$ws.Range("N83").Value = -2022283
$ws.Range("H97").Value = 662357.75
$ws.Range("I97").Value = 701196.5
$ws.Range("K97").Value = 701196.5
$ws.Range("M97").Value = -700700.5
$ws.Range("H107").Value = 1146.091
$ws.Range("I107").Value = 1227.4706
$ws.Range("K107").Value = 1227.4706
$ws.Range("M107").Value = 692.5293999999999
$ws.Range("H109").Value = 43817
$ws.Range("J109").Value = 48521.25
$ws.Range("L109").Value = 48521.25
$ws.Range("N109").Value = -50601.25
$ws.Range("H113").Value = 18521034
$ws.Range("I113").Value = 18521034
$ws.Range("K113").Value = 18521034
$ws.Range("M113").Value = -18518864
$ws.Range("H116").Value = 120461.336
$ws.Range("J116").Value = 120461.336
$ws.Range("L116").Value = 120461.336
$ws.Range("N116").Value = -129639.336
$ws.Range("H117").Value = 32666.334
$ws.Range("J117").Value = 32666.334
$ws.Range("L117").Value = 32666.334
$ws.Range("N117").Value = -39550.334
$ws.Range("H118").Value = 50000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 50000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 50000
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -53314
$ws.Range("H119").Value = 99999
$ws.Range("J119").Value = 99999
$ws.Range("L119").Value = 99999
$ws.Range("N119").Value = -109675
$ws.Range("H120").Value = 35918.145
$ws.Range("J120").Value = 35918.145
$ws.Range("L120").Value = 35918.145
$ws.Range("N120").Value = -45594.145
$ws.Range("H121").Value = 54895.668
$ws.Range("J121").Value = 54895.668
$ws.Range("L121").Value = 54895.668
$ws.Range("N121").Value = -58389.668
$ws.Range("H122").Value = 253368.73
$ws.Range("I122").Value = 317351.16
$ws.Range("J122").Value = 49788.273
$ws.Range("K122").Value = 952053.48
$ws.Range("L122").Value = 149364.819
$ws.Range("M122").Value = -949603.48
$ws.Range("N122").Value = -154264.819
$ws.Range("H130").Value = 47425.75
$ws.Range("J130").Value = 47401
$ws.Range("L130").Value = 47401
$ws.Range("N130").Value = -57441

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 960.86487
$ws.Range("I16").Value = 734.6061
$ws.Range("K16").Value = 734.6061
$ws.Range("M16").Value = -564.6061
$ws.Range("H63").Value = 75000
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 75000
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H82").Value = 70371720
$ws.Range("I82").Value = 95960984
$ws.Range("J82").Value = 1243
$ws.Range("K82").Value = 95960984
$ws.Range("L82").Value = 1243
$ws.Range("M82").Value = -95960623
$ws.Range("N82").Value = -1965
$ws.Range("H85").Value = 70371720
$ws.Range("I85").Value = 95960984
$ws.Range("J85").Value = 1243
$ws.Range("K85").Value = 95960984
$ws.Range("L85").Value = 1243
$ws.Range("M85").Value = -95959736
$ws.Range("N85").Value = -3739
$ws.Range("H93").Value = 15874269
$ws.Range("I93").Value = 20834348
$ws.Range("K93").Value = 20834348
$ws.Range("M93").Value = -20833100
$ws.Range("H132").Value = 8068.846
$ws.Range("I132").Value = 8434.348
$ws.Range("K132").Value = 25303.044
$ws.Range("M132").Value = -22773.044
$ws.Range("H134").Value = 77158.11
$ws.Range("J134").Value = 77158.11
$ws.Range("L134").Value = 77158.11
$ws.Range("N134").Value = -87298.11
$ws.Range("H136").Value = 37588.13
$ws.Range("I136").Value = 62288.766
$ws.Range("K136").Value = 186866.298
$ws.Range("M136").Value = -184316.298

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5000
$ws.Range("J2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("N2").Value = -5224
$ws.Range("H11").Value = 10005000
$ws.Range("I11").Value = 10005000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 10005000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -10004858
$ws.Range("N11").ClearContents()
$ws.Range("H18").Value = 9999
$ws.Range("I18").Value = 9999
$ws.Range("K18").Value = 9999
$ws.Range("M18").Value = -9826
$ws.Range("H32").Value = 5000
$ws.Range("I32").Value = 5000
$ws.Range("K32").Value = 5000
$ws.Range("M32").Value = -4683
$ws.Range("H51").Value = 16379.6
$ws.Range("I51").Value = 13999.5
$ws.Range("K51").Value = 13999.5
$ws.Range("M51").Value = -13489.5
$ws.Range("H81").Value = 18521332
$ws.Range("I81").Value = 18521332
$ws.Range("K81").Value = 37042664
$ws.Range("M81").Value = -37041603
$ws.Range("H84").Value = 18521332
$ws.Range("I84").Value = 18521332
$ws.Range("K84").Value = 185213320
$ws.Range("M84").Value = -185208016
$ws.Range("H96").Value = 7266.1113
$ws.Range("J96").Value = 8597.799999999999
$ws.Range("L96").Value = 8597.799999999999
$ws.Range("N96").Value = -11343.8
$ws.Range("H130").Value = 40090
$ws.Range("I130").Value = 40090
$ws.Range("K130").Value = 40090
$ws.Range("M130").Value = -35070
$ws.Range("H132").Value = 17259692
$ws.Range("I132").Value = 19610888
$ws.Range("K132").Value = 58832664
$ws.Range("M132").Value = -58830134
$ws.Range("H135").Value = 116560.55
$ws.Range("J135").Value = 116560.55
$ws.Range("L135").Value = 116560.55
$ws.Range("N135").Value = -126700.55
$ws.Range("H136").Value = 2203.721
$ws.Range("I136").Value = 2128.2927
$ws.Range("K136").Value = 6384.8781
$ws.Range("M136").Value = -3834.8781
